# Add more keywords/features to group "ห้า" (tag #5), which pushes the
# existing "หก" (tag #6) rows down by 3 rows (they had reused some of the
# blank spacer rows that used to separate the two groups).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: relocate the existing "หก" rows (96-103) down by three rows,
# to (99-106). Walk bottom-up so we never overwrite a row before it has
# been read.
for ($r = 103; $r -ge 96; $r--) {
    $srcA = $ws.Cells.Item($r, 1).Value2
    $srcB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 3, 1).Value = $srcA
    $ws.Cells.Item($r + 3, 2).Value = $srcB
}

# Step 2: fill the freed-up rows (93-97) with the five new keyword rows
# for group "ห้า", and clear what is now the spacer row (98).
$ws.Cells.Item(93, 1).Value = "ห้า"
$ws.Cells.Item(93, 2).Value = "ร้านที่ต้องจดทะเบียนอิเล็กทรอนิกส์"

$ws.Cells.Item(94, 1).Value = "ห้า "
$ws.Cells.Item(94, 2).Value = "ลักษณะร้านที่ต้องจดทะเบียนอิเล็กทรอนิกซ์"

$ws.Cells.Item(95, 1).Value = "ห้า"
$ws.Cells.Item(95, 2).Value = "ร้านประเภทไหนต้องจดทะเบียนอิเล็กทรอนิกซ์"

$ws.Cells.Item(96, 1).Value = "ห้า"
$ws.Cells.Item(96, 2).Value = "เปิดร้านใน shopee ต้องมีทะเบียนอิเล็กทรอนิกส์ไหม"

$ws.Cells.Item(97, 1).Value = "ห้า"
$ws.Cells.Item(97, 2).Value = "เปิดร้านใน lazada ต้องมีทะเบียนอิเล็กทรอนิกส์ไหม"

$ws.Range("A98:B98").ClearContents()

# Step 3: restore the selection to reflect where editing left off.
$ws.Range("B97").Select() | Out-Null
